$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test")

$values = @(
    'option1',
    'Тестовая запись №1',
    'Тестовая запись №2',
    'Тестовая запись №3',
    'Тестовая запись №4',
    'Тестовая запись №5',
    'Тестовая запись №6',
    'Тестовая запись №7',
    'Тестовая запись №8',
    'Тестовая запись №9',
    'Тестовая запись №10',
    'Тестовая запись №11',
    'Тестовая запись №12',
    'Тестовая запись №13',
    'Тестовая запись №14',
    'Тестовая запись №15',
    'Тестовая запись №16',
    'Тестовая запись №17',
    'Тестовая запись №18',
    'Тестовая запись №19',
    'Тестовая запись №20',
    'Acer',
    'Toyota',
    'HP',
    'HP',
    'Toyota',
    'ASUS',
    'HP',
    'Intel',
    'HP',
    'Ford',
    'Ford',
    'Intel',
    'Samsung',
    'Ford',
    'Toyota',
    'AMD',
    'HP',
    'Acer',
    'Toyota',
    'Toyota',
    'Intel',
    'Ford',
    'Acer',
    'Samsung',
    'Samsung',
    'Ford',
    'ASUS',
    'Samsung',
    'Intel',
    'Acer',
    'Ford',
    'Samsung',
    'HP',
    'HP',
    'ASUS',
    'ASUS',
    'Ford',
    'Acer',
    'Samsung',
    'HP',
    'Toyota',
    'AMD',
    'Acer',
    'ASUS',
    'Acer',
    'ASUS',
    'ASUS',
    'Acer',
    'HP',
    'Samsung',
    'AMD',
    'Samsung',
    'Toyota',
    'HP',
    'Samsung',
    'Ford',
    'HP',
    'Ford',
    'Samsung',
    'AMD',
    'HP',
    'AMD',
    'Samsung',
    'AMD',
    'AMD',
    'Samsung',
    'Acer',
    'Acer',
    'Intel',
    'Samsung',
    'Acer',
    'HP',
    'HP',
    'ASUS',
    'Toyota',
    'Samsung',
    'Samsung',
    'Samsung',
    'Samsung',
    'Acer',
    'Intel',
    'Intel',
    'Acer',
    'Samsung',
    'Ford',
    'Acer',
    'Samsung',
    'AMD',
    'Toyota',
    'HP',
    'Ford',
    'Samsung',
    'Acer',
    'Ford',
    'Toyota',
    'Toyota',
    'Acer',
    'ASUS',
    'Samsung',
    'Acer'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

